$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Row 7 holds the "Experimental" property; it previously had no value in
# column B. FHIR requires the "experimental" element, so add it as the
# text "true" (matches the sheet's existing text-based boolean convention,
# e.g. row 15 "Immutable" -> "BooleanType[null]").
#
# A leading apostrophe forces the literal to be stored as text rather than
# being auto-coerced into a native Excel boolean; copying the neighboring
# "Property" cell's format back over it afterwards clears the resulting
# quote-prefix flag and keeps the cell on the sheet's normal "value" style
# instead of minting a new one.
$wsMeta.Range("B7").Value = "'true"
$wsMeta.Range("A7").Copy() | Out-Null
$wsMeta.Range("B7").PasteSpecial(-4122) | Out-Null

# Row 8 holds the "Date" property; bump it to the new export timestamp.
$wsMeta.Range("B8").Value = "2023-02-01T09:05:11-06:00"

$excel.CutCopyMode = 0
